$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the blank placeholder row (row 3, containing "   "), which
# shifts "contact_info" (and everything below it) up by one row.
$ws.Rows(3).Delete()

# Update the selection to reflect where the edit took place.
$ws.Range("A3").Select()
